# Added okadaic acid to the analysis:
# Append 5 new data rows (17-21) to the docking-analysis results table on
# the active sheet, then move the selection the way the author left it
# (single cell P21) and nudge the viewport to the right (column E onward)
# to mirror the saved sheetView in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write each new row left-to-right (A -> O) so Excel records the row's
# column span correctly.

# Row 17
$ws.Cells.Item(17, 1).Value = -59
$ws.Cells.Item(17, 2).Value = 12
$ws.Cells.Item(17, 3).Value = 0.7
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = -24.4
$ws.Cells.Item(17, 6).Value = 13.2
$ws.Cells.Item(17, 7).Value = -335.3
$ws.Cells.Item(17, 8).Value = 26.1
$ws.Cells.Item(17, 9).Value = -13.4
$ws.Cells.Item(17, 10).Value = 1.9
$ws.Cells.Item(17, 11).Value = 123.5
$ws.Cells.Item(17, 12).Value = 35.75
$ws.Cells.Item(17, 13).Value = 1234.8
$ws.Cells.Item(17, 14).Value = 11.5
$ws.Cells.Item(17, 15).Value = -1.6

# Row 18
$ws.Cells.Item(18, 1).Value = -57.5
$ws.Cells.Item(18, 2).Value = 9.5
$ws.Cells.Item(18, 3).Value = 0.7
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = -22.7
$ws.Cells.Item(18, 6).Value = 10.3
$ws.Cells.Item(18, 7).Value = -338.3
$ws.Cells.Item(18, 8).Value = 22.8
$ws.Cells.Item(18, 9).Value = -13.3
$ws.Cells.Item(18, 10).Value = 1.8
$ws.Cells.Item(18, 11).Value = 123.1
$ws.Cells.Item(18, 12).Value = 35.59
$ws.Cells.Item(18, 13).Value = 1237.6
$ws.Cells.Item(18, 14).Value = 11.7
$ws.Cells.Item(18, 15).Value = -1.7

# Row 19
$ws.Cells.Item(19, 1).Value = -64.3
$ws.Cells.Item(19, 2).Value = 9.7
$ws.Cells.Item(19, 3).Value = 0.2
$ws.Cells.Item(19, 4).Value = 0.1
$ws.Cells.Item(19, 5).Value = -24.7
$ws.Cells.Item(19, 6).Value = 6.3
$ws.Cells.Item(19, 7).Value = -392.5
$ws.Cells.Item(19, 8).Value = 48.5
$ws.Cells.Item(19, 9).Value = -12.6
$ws.Cells.Item(19, 10).Value = 0.3
$ws.Cells.Item(19, 11).Value = 121.7
$ws.Cells.Item(19, 12).Value = 42.94
$ws.Cells.Item(19, 13).Value = 1205.5
$ws.Cells.Item(19, 14).Value = 8.2
$ws.Cells.Item(19, 15).Value = -1.4

# Row 20
$ws.Cells.Item(20, 1).Value = -65.8
$ws.Cells.Item(20, 2).Value = 9.3
$ws.Cells.Item(20, 3).Value = 0.7
$ws.Cells.Item(20, 4).Value = 0.1
$ws.Cells.Item(20, 5).Value = -29.7
$ws.Cells.Item(20, 6).Value = 4.4
$ws.Cells.Item(20, 7).Value = -363.8
$ws.Cells.Item(20, 8).Value = 45.5
$ws.Cells.Item(20, 9).Value = -12.8
$ws.Cells.Item(20, 10).Value = 0.8
$ws.Cells.Item(20, 11).Value = 130.1
$ws.Cells.Item(20, 12).Value = 30.73
$ws.Cells.Item(20, 13).Value = 1212
$ws.Cells.Item(20, 14).Value = 5
$ws.Cells.Item(20, 15).Value = -1.9

# Row 21
$ws.Cells.Item(21, 1).Value = -72.9
$ws.Cells.Item(21, 2).Value = 10.9
$ws.Cells.Item(21, 3).Value = 0.7
$ws.Cells.Item(21, 4).Value = 0.1
$ws.Cells.Item(21, 5).Value = -34.4
$ws.Cells.Item(21, 6).Value = 11.3
$ws.Cells.Item(21, 7).Value = -340.5
$ws.Cells.Item(21, 8).Value = 39.6
$ws.Cells.Item(21, 9).Value = -14.8
$ws.Cells.Item(21, 10).Value = 1.4
$ws.Cells.Item(21, 11).Value = 104.6
$ws.Cells.Item(21, 12).Value = 21.86
$ws.Cells.Item(21, 13).Value = 1224
$ws.Cells.Item(21, 14).Value = 6.9
$ws.Cells.Item(21, 15).Value = -2.4

# Scroll the viewport right so column E is the left-most visible column
# (mirrors topLeftCell="E1" in the saved view) and leave the selection on
# P21, matching the author's final cursor position.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("P21").Select()
